# Auto-generated data-driven update of cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='26.702.84'; E='  +1.57%  ' },
    @{ Row=3; D='1.636.60'; E='  +1.65%  ' },
    @{ Row=4; E='  -0.01%  ' },
    @{ Row=5; D='213.63'; E='  +0.30%  ' },
    @{ Row=6; E='  +0.01%  ' },
    @{ Row=7; E='  +1.12%  ' },
    @{ Row=8; E='  +0.71%  ' },
    @{ Row=9; D='0.0620'; E='  +0.84%  ' },
    @{ Row=10; D='19.07'; E='  +3.66%  ' },
    @{ Row=11; E='  +2.29%  ' },
    @{ Row=12; D='1.865.22'; E='  +1.74%  ' },
    @{ Row=13; D='1.619.55'; E='  +0.73%  ' },
    @{ Row=14; E='  +0.32%  ' },
    @{ Row=15; E='  +1.74%  ' },
    @{ Row=16; D='26.680.85'; E='  +1.45%  ' },
    @{ Row=17; D='63.28'; E='  +2.35%  ' },
    @{ Row=18; E='  +0.68%  ' },
    @{ Row=19; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.00'; E='  -0.02%  ' },
    @{ Row=20; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='208.33'; E='  +2.41%  ' },
    @{ Row=21; E='  +0.76%  ' },
    @{ Row=22; D='9.39'; E='  +0.89%  ' },
    @{ Row=23; E='  +1.42%  ' },
    @{ Row=24; E='  -0.79%  ' },
    @{ Row=25; D='145.76'; E='  +1.01%  ' },
    @{ Row=26; E='  +0.02%  ' },
    @{ Row=27; D='0.120'; E='  -1.58%  ' },
    @{ Row=28; E='  +1.01%  ' },
    @{ Row=29; D='6.68'; E='  +1.81%  ' },
    @{ Row=30; D='0.0517'; E='  +5.62%  ' },
    @{ Row=31; D='1.18'; E='  +0.45%  ' },
    @{ Row=32; E='  +0.86%  ' },
    @{ Row=33; D='2.95'; E='  +0.16%  ' },
    @{ Row=34; E='  +1.31%  ' },
    @{ Row=35; E='  -0.49%  ' },
    @{ Row=36; D='1.167.00'; E='  +0.67%  ' },
    @{ Row=37; E='  +0.35%  ' },
    @{ Row=38; D='0.814'; E='  +2.06%  ' },
    @{ Row=40; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.504'; E='  +0.31%  ' },
    @{ Row=41; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.32'; E='  -0.12%  ' },
    @{ Row=42; D='5.41'; E='  +2.96%  ' },
    @{ Row=43; D='0.795'; E='  +1.23%  ' },
    @{ Row=44; D='1.773.89'; E='  +1.57%  ' },
    @{ Row=45; D='92.45'; E='  +0.81%  ' },
    @{ Row=46; E='  +0.81%  ' },
    @{ Row=47; D='54.70'; E='  +0.56%  ' },
    @{ Row=48; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.0511'; E='  +0.86%  ' },
    @{ Row=49; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='7.62'; E='  +4.73%  ' },
    @{ Row=50; B='Mantle'; C='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D='0.410'; E='  +0.82%  ' },
    @{ Row=51; B='USDD'; C='https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'; D='1.00'; E='  -0.10%  ' }
)

$colMap = @{ 'B' = 2; 'C' = 3; 'D' = 4; 'E' = 5 }

foreach ($update in $updates) {
    $row = $update['Row']
    foreach ($col in @('B','C','D','E')) {
        if ($update.ContainsKey($col)) {
            $value = $update[$col]
            $cell = $ws.Cells.Item($row, $colMap[$col])
            if ($col -eq 'D') {
                # Force text interpretation so numeric-looking strings (e.g. "1.00", "213.63")
                # keep their exact original formatting instead of being parsed as numbers.
                $cell.NumberFormat = "@"
            }
            $cell.Value = $value
        }
    }
}
